$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.849458694458008
$ws.Range("B1").Value = 3.187530279159546
$ws.Range("C1").Value = 2.93264102935791
$ws.Range("D1").Value = 3.512891054153442
$ws.Range("E1").Value = 3.891805648803711
